$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 336, shifting existing rows 336-397 down to 337-398
$ws.Rows.Item(336).Insert()

# Populate the newly inserted row 336 with data
$ws.Cells.Item(336, 1).Value = 8
$ws.Cells.Item(336, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(336, 3).Value = "Coquimbo"
$ws.Cells.Item(336, 4).Value = 45015
$ws.Cells.Item(336, 5).Value = 4
$ws.Cells.Item(336, 6).Value = 100112012
$ws.Cells.Item(336, 7).Value = "Espinaca"
$ws.Cells.Item(336, 8).Value = "Sin especificar"
$ws.Cells.Item(336, 9).Value = "Primera"
$ws.Cells.Item(336, 10).Value = 1500
$ws.Cells.Item(336, 11).Value = 450
$ws.Cells.Item(336, 12).Value = 500
$ws.Cells.Item(336, 13).Value = 475
$ws.Cells.Item(336, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(336, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(336, 16).Value = 950
$ws.Cells.Item(336, 17).Value = 0.5
$ws.Cells.Item(336, 18).Value = "Hortaliza"
